$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PercentText {
    param($addr, $text)
    $c = $ws.Range($addr)
    # Leading apostrophe forces literal text so Excel doesn't auto-convert
    # the "NN.N%" string into a numeric percentage (which would change the
    # cell's stored type/value and saved style).
    $c.Value = "'" + $text
    # Re-applying the "Normal" style clears the sticky percent-number-format
    # that Excel's smart entry would otherwise leave behind, then we restore
    # the original center/center alignment used throughout this sheet.
    $c.Style = "Normal"
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
}

# --- Class Statistics (K2:L10 block) ---
$ws.Range("L6").Value = 114
$ws.Range("L8").Value = 114
Set-PercentText "L9" "48.7%"
Set-PercentText "L10" "78.0%"

# --- Group Statistics rows (O/P/Q/R/S columns) ---
$ws.Range("O16").Value = 10
$ws.Range("Q16").Value = 10
Set-PercentText "R16" "50.0%"
Set-PercentText "S16" "71.9%"

$ws.Range("O17").Value = 10
$ws.Range("Q17").Value = 10
Set-PercentText "R17" "50.0%"
Set-PercentText "S17" "60.0%"

$ws.Range("O18").Value = 10
$ws.Range("Q18").Value = 10
Set-PercentText "R18" "50.0%"
Set-PercentText "S18" "82.9%"

$ws.Range("O24").Value = 10
$ws.Range("Q24").Value = 10
Set-PercentText "R24" "50.0%"
Set-PercentText "S24" "71.9%"

$ws.Range("O25").Value = 10
$ws.Range("Q25").Value = 10
Set-PercentText "R25" "50.0%"
Set-PercentText "S25" "74.8%"

$ws.Range("O26").Value = 10
$ws.Range("Q26").Value = 10
Set-PercentText "R26" "50.0%"
Set-PercentText "S26" "71.7%"

# --- Sessions that moved from "Pending" to "Recorded" ---
# Each of these rows switches its row style from the yellow "Pending" fill
# (style index 6) to the green "Recorded" fill (style index 2), gains a
# "Recorded By" e-mail, and its attendance count changes from "0/N" to the
# real "<attended>/N".
function Set-SessionRecorded {
    param($row, $attended, $total)

    $range = $ws.Range("A" + $row + ":I" + $row)
    # Copy the look of an already-"Recorded" row (style 2, green fill,
    # centered) onto this row.
    $range.Style = $ws.Range("A24:I24").Style

    $ws.Range("G" + $row).Value = "dnasr281@gmail.com"
    $ws.Range("H" + $row).Value = $attended + "/" + $total
    $ws.Range("I" + $row).Value = "Recorded"
}

Set-SessionRecorded 27  "27" "31"
Set-SessionRecorded 47  "13" "18"
Set-SessionRecorded 67  "15" "21"
Set-SessionRecorded 182 "22" "27"
Set-SessionRecorded 202 "21" "29"
Set-SessionRecorded 222 "19" "29"
